$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHAM_CONG_HE_THONG")

# Rows 2-14 -> 17:23, rows 15-20 -> 17:24 (all date 2024-07-17)
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("D$r").Value = "2024-07-17T17:23:00.000Z"
}
for ($r = 15; $r -le 20; $r++) {
    $ws.Range("D$r").Value = "2024-07-17T17:24:00.000Z"
}
